$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New logged-hours rows for the week of 1/31-2/11/2018, plus "really ___" notes
# in the new column F explaining the true day worked for entries that were
# logged a bit late. Values/order below mirror the exact sequence the data
# was entered in, so new shared strings land at the right table positions.

$ws.Cells.Item(290, 1).Value = "Friday, Feb 2, 2018"
$ws.Cells.Item(290, 2).Value = 0.58333333333333337
$ws.Cells.Item(290, 3).Value = 0.61458333333333337

$ws.Cells.Item(291, 1).Value = "Friday, Feb 2, 2018"
$ws.Cells.Item(291, 2).Value = 0.83333333333333337
$ws.Cells.Item(291, 3).Value = 0.98958333333333337

$ws.Cells.Item(292, 1).Value = "Thursday, Feb 8, 2018"
$ws.Cells.Item(292, 2).Value = 0.70833333333333337
$ws.Cells.Item(292, 3).Value = 0.76041666666666663

$ws.Cells.Item(289, 6).Value = "really Wednesday 1/31"

$ws.Cells.Item(291, 6).Value = "really Thursday 2/1"

$ws.Cells.Item(292, 6).Value = "really Friday 2/2"

$ws.Cells.Item(289, 1).Copy() | Out-Null
$ws.Cells.Item(293, 6).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(293, 6).Value = "really Saturday 2/3"

$ws.Cells.Item(293, 1).Value = "Friday, Feb 9, 2018"
$ws.Cells.Item(293, 2).Value = 0.69791666666666663
$ws.Cells.Item(293, 3).Value = 0.79166666666666663

$ws.Cells.Item(294, 6).Value = "really Sunday 2/4"

$ws.Cells.Item(295, 6).Value = "really Thursday 2/8"

$ws.Cells.Item(295, 1).Value = "Saturday, Feb 10, 2018"
$ws.Cells.Item(295, 2).Value = 0.625
$ws.Cells.Item(295, 3).Value = 0.64583333333333337

$ws.Cells.Item(296, 1).Value = "Saturday, Feb 10, 2018"
$ws.Cells.Item(296, 2).Value = 0.91666666666666663
$ws.Cells.Item(296, 3).Value = 0.95833333333333337
$ws.Cells.Item(296, 6).Value = "really Thursday 2/8"

$ws.Cells.Item(297, 1).Value = "Sunday, Feb 11, 2018"
$ws.Cells.Item(297, 2).Value = 0.58333333333333337
$ws.Cells.Item(297, 3).Value = 0.6875
$ws.Cells.Item(297, 6).Value = "really Friday 2/9"

# Give the new notes column a sensible width
$ws.Columns.Item(6).ColumnWidth = 21

# Restore the view: scrolled down a bit further, with D297 now selected
$win = $excel.ActiveWindow
$win.ScrollRow = 285
$win.ScrollColumn = 1
$ws.Range("D297").Select()
